# Update the "Förändrad" (Changed) date column (C) for all data rows (2-471)
# from serial 45171 (2023-09-02) to serial 45172 (2023-09-03).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 471; $r++) {
    $ws.Cells.Item($r, 3).Value = 45172
}

# Rows 470 and 471 had their "Beteckning" (A) and "Area (ha)" (G) values swapped.
$ws.Cells.Item(470, 1).Value = "A 40725-2023"
$ws.Cells.Item(470, 7).Value = 1.4
$ws.Cells.Item(471, 1).Value = "A 40731-2023"
$ws.Cells.Item(471, 7).Value = 1.8
